$d = $word.ActiveDocument

function Find-ParagraphContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# 1. Remove the existing _GoBack bookmark - it sits inside the paragraph
#    that is about to be deleted ("Generation works for properties...").
#    It will be re-created at its new location (inside the "invalid part
#    is underlined" paragraph) in step 2.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2. Split the run "valid part is underlined and the OK button is
#    disabled." into "valid part is un" + bookmark + "derlined and the OK
#    button is disabled." by locating the split point and inserting a
#    zero-length bookmark named _GoBack there.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "valid part is un"
$find.Execute() | Out-Null
$splitPoint = $find.Parent.Duplicate
$splitPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# 3. Remove the paragraphs between "If the Destination Package does not
#    exist it is created." and "What happens if you try to generate using
#    a superclass..." - i.e. the "What methods are generated?" through
#    "Check generator works for matchers in the default package." block.
$block1Start = Find-ParagraphContaining "What methods are generated"
$block1End = Find-ParagraphContaining "Check generator works for matchers in the default package"
$d.Range($block1Start.Range.Start, $block1End.Range.End).Delete()

# 4. Remove everything after "What happens if you try to generate using a
#    superclass..." through to the end of the document ("Generated class
#    should be public" .. "Generation works for properties with primitive
#    and non-primitive types.").
$keepLast = Find-ParagraphContaining "What happens if you try to generate"
$block2Start = $keepLast.Next()
$block2End = $d.Paragraphs.Last
$d.Range($block2Start.Range.Start, $block2End.Range.End).Delete()
